# Scheduled-runner update: refresh cached market-board figures
# (currentAveragePrice / NQ / HQ / LevePrice* / LeveProfit*) for a batch
# of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 416.55554
$ws.Range("I2").Value = 289.8
$ws.Range("K2").Value = 289.8
$ws.Range("M2").Value = -176.8
$ws.Range("H33").Value = 247.76923
$ws.Range("I33").Value = 127.85714
$ws.Range("K33").Value = 127.85714
$ws.Range("M33").Value = 101.14286
$ws.Range("H40").Value = 3275.125
$ws.Range("I40").Value = 5799.5
$ws.Range("K40").Value = 5799.5
$ws.Range("M40").Value = -5624.5
$ws.Range("H137").Value = 1090.3684
$ws.Range("I137").Value = 1095.3889
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 3286.1667
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -736.1666999999998
$ws.Range("N137").Value = -8100
$ws.Range("H138").Value = 1497.42
$ws.Range("I138").Value = 682.9729599999999
$ws.Range("J138").Value = 1975.746
$ws.Range("K138").Value = 2048.91888
$ws.Range("L138").Value = 5927.238
$ws.Range("M138").Value = 3091.08112
$ws.Range("N138").Value = -16207.238

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 23000
$ws.Range("J96").Value = 23000
$ws.Range("L96").Value = 23000
$ws.Range("N96").Value = -28492
$ws.Range("H122").Value = 1970.6666
$ws.Range("I122").Value = 1956
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5868
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3418
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 3526.25
$ws.Range("I132").Value = 3238.4546
$ws.Range("K132").Value = 9715.363799999999
$ws.Range("M132").Value = -7185.363799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1754.9131
$ws.Range("I107").Value = 1253.8462
$ws.Range("J107").Value = 2406.3
$ws.Range("K107").Value = 1253.8462
$ws.Range("L107").Value = 2406.3
$ws.Range("M107").Value = 666.1538
$ws.Range("N107").Value = -6246.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1075.5714
$ws.Range("I31").Value = 973.6923
$ws.Range("K31").Value = 973.6923
$ws.Range("M31").Value = -678.6923
$ws.Range("H34").Value = 1075.5714
$ws.Range("I34").Value = 973.6923
$ws.Range("K34").Value = 973.6923
$ws.Range("M34").Value = -771.6923
$ws.Range("H50").Value = 14940.429
$ws.Range("I50").Value = 2791.5
$ws.Range("K50").Value = 2791.5
$ws.Range("M50").Value = -2166.5
$ws.Range("H99").Value = 2393876.5
$ws.Range("I99").Value = 3290856.5
$ws.Range("J99").Value = 1930
$ws.Range("K99").Value = 3290856.5
$ws.Range("L99").Value = 1930
$ws.Range("M99").Value = -3289358.5
$ws.Range("N99").Value = -4926
$ws.Range("H122").Value = 1206
$ws.Range("I122").Value = 1012
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 3036
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -586
$ws.Range("N122").Value = -9100
$ws.Range("H126").Value = 2393876.5
$ws.Range("I126").Value = 3290856.5
$ws.Range("J126").Value = 1930
$ws.Range("K126").Value = 9872569.5
$ws.Range("L126").Value = 5790
$ws.Range("M126").Value = -9870099.5
$ws.Range("N126").Value = -10730
$ws.Range("H132").Value = 9184.25
$ws.Range("I132").Value = 13400.111
$ws.Range("J132").Value = 3763.8572
$ws.Range("K132").Value = 40200.333
$ws.Range("L132").Value = 11291.5716
$ws.Range("M132").Value = -37670.333
$ws.Range("N132").Value = -16351.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 111.1
$ws.Range("J2").Value = 203
$ws.Range("L2").Value = 1218
$ws.Range("N2").Value = -1444
$ws.Range("H5").Value = 1408.591
$ws.Range("I5").Value = 1408.591
$ws.Range("K5").Value = 4225.772999999999
$ws.Range("M5").Value = -4113.772999999999
$ws.Range("H32").Value = 2218.182
$ws.Range("J32").Value = 2218.182
$ws.Range("L32").Value = 6654.545999999999
$ws.Range("N32").Value = -7220.545999999999
$ws.Range("H131").Value = 13514751
$ws.Range("J131").Value = 1340.4546
$ws.Range("L131").Value = 4021.3638
$ws.Range("N131").Value = -14101.3638
$ws.Range("H135").Value = 1408.591
$ws.Range("I135").Value = 1408.591
$ws.Range("K135").Value = 12677.319
$ws.Range("M135").Value = -10142.319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.818184
$ws.Range("I2").Value = 98.166664
$ws.Range("J2").Value = 119.4
$ws.Range("K2").Value = 98.166664
$ws.Range("L2").Value = 119.4
$ws.Range("M2").Value = 14.833336
$ws.Range("N2").Value = -345.4
$ws.Range("H12").Value = 4454739
$ws.Range("I12").Value = 4212333.5
$ws.Range("K12").Value = 4212333.5
$ws.Range("M12").Value = -4212193.5
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H43").Value = 6226
$ws.Range("I43").Value = 1200
$ws.Range("J43").Value = 6682.909
$ws.Range("K43").Value = 1200
$ws.Range("L43").Value = 6682.909
$ws.Range("M43").Value = -1049
$ws.Range("N43").Value = -6984.909
$ws.Range("H126").Value = 2218.6667
$ws.Range("I126").Value = 1812.1111
$ws.Range("J126").Value = 2523.5833
$ws.Range("K126").Value = 5436.3333
$ws.Range("L126").Value = 7570.749899999999
$ws.Range("M126").Value = -2966.3333
$ws.Range("N126").Value = -12510.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1967.8334
$ws.Range("I7").Value = 1975.5
$ws.Range("K7").Value = 1975.5
$ws.Range("M7").Value = -1863.5
$ws.Range("H40").Value = 3039.8
$ws.Range("I40").Value = 2759.2
$ws.Range("J40").Value = 3601
$ws.Range("K40").Value = 2759.2
$ws.Range("L40").Value = 3601
$ws.Range("M40").Value = -2623.2
$ws.Range("N40").Value = -3873
$ws.Range("H46").Value = 2100
$ws.Range("I46").Value = 1933.3334
$ws.Range("J46").Value = 2200
$ws.Range("K46").Value = 1933.3334
$ws.Range("L46").Value = 2200
$ws.Range("M46").Value = -1745.3334
$ws.Range("N46").Value = -2576
$ws.Range("H122").Value = 28336118
$ws.Range("I122").Value = 70836160
$ws.Range("J122").Value = 2759
$ws.Range("K122").Value = 212508480
$ws.Range("L122").Value = 8277
$ws.Range("M122").Value = -212506030
$ws.Range("N122").Value = -13177
$ws.Range("H126").Value = 1967.8334
$ws.Range("I126").Value = 1975.5
$ws.Range("K126").Value = 5926.5
$ws.Range("M126").Value = -3456.5
$ws.Range("H140").Value = 50776.668
$ws.Range("J140").Value = 50776.668
$ws.Range("L140").Value = 50776.668
$ws.Range("N140").Value = -61136.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4424.095
$ws.Range("I81").Value = 648.5
$ws.Range("K81").Value = 1297
$ws.Range("M81").Value = -236
$ws.Range("H84").Value = 4424.095
$ws.Range("I84").Value = 648.5
$ws.Range("K84").Value = 6485
$ws.Range("M84").Value = -1181
$ws.Range("H122").Value = 15296784
$ws.Range("I122").Value = 15296784
$ws.Range("K122").Value = 45890352
$ws.Range("M122").Value = -45887902
$ws.Range("H126").Value = 222223420
$ws.Range("I126").Value = 222223420
$ws.Range("K126").Value = 666670260
$ws.Range("M126").Value = -666667790
